$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply each cell update as a literal text value, preserving the original
# "General" number format / default style once the text has been written.
$c = $ws.Range('D2')
$c.NumberFormat = '@'
$c.Value = '68.294.96'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E2')
$c.NumberFormat = '@'
$c.Value = '  +0.57%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D3')
$c.NumberFormat = '@'
$c.Value = '2.642.26'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E3')
$c.NumberFormat = '@'
$c.Value = '  +0.46%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E4')
$c.NumberFormat = '@'
$c.Value = '  +0.07%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D5')
$c.NumberFormat = '@'
$c.Value = '598.71'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E5')
$c.NumberFormat = '@'
$c.Value = '  +0.15%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D6')
$c.NumberFormat = '@'
$c.Value = '154.64'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E6')
$c.NumberFormat = '@'
$c.Value = '  +0.60%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E7')
$c.NumberFormat = '@'
$c.Value = '  +0.03%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E8')
$c.NumberFormat = '@'
$c.Value = '  -0.68%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D9')
$c.NumberFormat = '@'
$c.Value = '2.642.82'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E9')
$c.NumberFormat = '@'
$c.Value = '  +0.51%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E10')
$c.NumberFormat = '@'
$c.Value = '  +8.11%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E11')
$c.NumberFormat = '@'
$c.Value = '  -0.62%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D12')
$c.NumberFormat = '@'
$c.Value = '5.26'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E12')
$c.NumberFormat = '@'
$c.Value = '  +1.00%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D13')
$c.NumberFormat = '@'
$c.Value = '0.354'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E13')
$c.NumberFormat = '@'
$c.Value = '  +2.21%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E14')
$c.NumberFormat = '@'
$c.Value = '  +2.81%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D15')
$c.NumberFormat = '@'
$c.Value = '27.97'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E15')
$c.NumberFormat = '@'
$c.Value = '  +1.42%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D16')
$c.NumberFormat = '@'
$c.Value = '3.123.55'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E16')
$c.NumberFormat = '@'
$c.Value = '  +0.59%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D17')
$c.NumberFormat = '@'
$c.Value = '68.189.73'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E17')
$c.NumberFormat = '@'
$c.Value = '  +0.60%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D18')
$c.NumberFormat = '@'
$c.Value = '2.657.80'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E18')
$c.NumberFormat = '@'
$c.Value = '  +1.36%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D19')
$c.NumberFormat = '@'
$c.Value = '11.36'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E19')
$c.NumberFormat = '@'
$c.Value = '  -0.78%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D20')
$c.NumberFormat = '@'
$c.Value = '364.22'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E20')
$c.NumberFormat = '@'
$c.Value = '  -1.81%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D21')
$c.NumberFormat = '@'
$c.Value = '7.46'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E21')
$c.NumberFormat = '@'
$c.Value = '  +0.27%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D22')
$c.NumberFormat = '@'
$c.Value = '4.37'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E22')
$c.NumberFormat = '@'
$c.Value = '  +3.34%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E23')
$c.NumberFormat = '@'
$c.Value = '  +1.97%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('B24')
$c.NumberFormat = '@'
$c.Value = 'Litecoin'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('C24')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D24')
$c.NumberFormat = '@'
$c.Value = '75.73'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E24')
$c.NumberFormat = '@'
$c.Value = '  +4.99%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('B25')
$c.NumberFormat = '@'
$c.Value = 'SuiNetwork'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('C25')
$c.NumberFormat = '@'
$c.Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D25')
$c.NumberFormat = '@'
$c.Value = '2.06'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E25')
$c.NumberFormat = '@'
$c.Value = '  -0.90%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D26')
$c.NumberFormat = '@'
$c.Value = '1.00'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E26')
$c.NumberFormat = '@'
$c.Value = '  +0.11%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D27')
$c.NumberFormat = '@'
$c.Value = '9.77'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E27')
$c.NumberFormat = '@'
$c.Value = '  -1.07%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E28')
$c.NumberFormat = '@'
$c.Value = '  +2.26%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D29')
$c.NumberFormat = '@'
$c.Value = '2.778.91'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E29')
$c.NumberFormat = '@'
$c.Value = '  +0.73%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E30')
$c.NumberFormat = '@'
$c.Value = '  -0.44%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D31')
$c.NumberFormat = '@'
$c.Value = '562.36'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E31')
$c.NumberFormat = '@'
$c.Value = '  -2.33%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D32')
$c.NumberFormat = '@'
$c.Value = '8.05'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E32')
$c.NumberFormat = '@'
$c.Value = '  +1.74%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D33')
$c.NumberFormat = '@'
$c.Value = '1.41'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E33')
$c.NumberFormat = '@'
$c.Value = '  +0.62%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D34')
$c.NumberFormat = '@'
$c.Value = '1.86'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E34')
$c.NumberFormat = '@'
$c.Value = '  +1.26%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E35')
$c.NumberFormat = '@'
$c.Value = '  +2.31%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E36')
$c.NumberFormat = '@'
$c.Value = '  +0.05%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E37')
$c.NumberFormat = '@'
$c.Value = '  +4.16%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D38')
$c.NumberFormat = '@'
$c.Value = '161.23'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E38')
$c.NumberFormat = '@'
$c.Value = '  +1.11%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D39')
$c.NumberFormat = '@'
$c.Value = '19.33'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E39')
$c.NumberFormat = '@'
$c.Value = '  +0.85%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D40')
$c.NumberFormat = '@'
$c.Value = '0.374'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E40')
$c.NumberFormat = '@'
$c.Value = '  +1.76%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E41')
$c.NumberFormat = '@'
$c.Value = '  -0.32%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D42')
$c.NumberFormat = '@'
$c.Value = '5.35'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E42')
$c.NumberFormat = '@'
$c.Value = '  -0.09%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D43')
$c.NumberFormat = '@'
$c.Value = '0.0₆0340'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E43')
$c.NumberFormat = '@'
$c.Value = '  +2.98%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E44')
$c.NumberFormat = '@'
$c.Value = '  +0.43%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E46')
$c.NumberFormat = '@'
$c.Value = '  +0.01%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D47')
$c.NumberFormat = '@'
$c.Value = '40.61'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E47')
$c.NumberFormat = '@'
$c.Value = '  +0.96%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D48')
$c.NumberFormat = '@'
$c.Value = '156.18'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E48')
$c.NumberFormat = '@'
$c.Value = '  +0.38%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E49')
$c.NumberFormat = '@'
$c.Value = '  +1.95%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D50')
$c.NumberFormat = '@'
$c.Value = '1.70'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E50')
$c.NumberFormat = '@'
$c.Value = '  +0.56%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('D51')
$c.NumberFormat = '@'
$c.Value = '21.82'
$c.NumberFormat = 'General'
$c.Style = 'Normal'

$c = $ws.Range('E51')
$c.NumberFormat = '@'
$c.Value = '  -0.67%  '
$c.NumberFormat = 'General'
$c.Style = 'Normal'

